# Fixed Bug with no updating group after choosing combobox
# Updates the attendance sheet: rename the "Посещения" header to
# "Кол-во посещений" and append attendance rows for three more students.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fix -----------------------------------------------------------
$ws.Range("D1").Value = "Кол-во посещений"

# --- Row 3: Ростислав, same group/date/count as row 2 ---------------------
$ws.Range("A2:D2").Copy()
$ws.Range("A3").PasteSpecial()
$ws.Range("A3").Value = "Ростислав"

# --- Row 4: Неважно кто, new group, same date/count as row 2 --------------
$ws.Range("A2:D2").Copy()
$ws.Range("A4").PasteSpecial()
$ws.Range("A4").Value = "Неважно кто"

# Write the new group number as text (not a number) without touching the
# cell style: render it through a formula, then paste back as a value.
$ws.Range("Z1").Formula = "=TEXT(10702323,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# --- Row 5: Баклажан, same group as row 4, same date/count as row 2 -------
$ws.Range("A4:D4").Copy()
$ws.Range("A5").PasteSpecial()
$ws.Range("A5").Value = "Баклажан"

$ws.Range("A1").Select()
